# Q1 2022 Fiscal Data update
# Adds Jan-Dec 2022 month columns (PR:QC) to the "ngcor" sheet (only
# Jan-Apr 2022 actually carry data, matching "Latest Data: April 2022"),
# and refreshes the "metadata" sheet's availability / latest-data labels.

$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------
# Sheet "ngcor": append 12 new month columns (PR..QC) for 2022
# ---------------------------------------------------------------------
$ws = $wb.Worksheets.Item("ngcor")

$monthCols = @("PR","PS","PT","PU","PV","PW","PX","PY","PZ","QA","QB","QC")
$monthNames = @("January 2022","February 2022","March 2022","April 2022","May 2022","June 2022","July 2022","August 2022","September 2022","October 2022","November 2022","December 2022")

for ($i = 0; $i -lt $monthCols.Length; $i++) {
    $col = $monthCols[$i]
    # Format header cell as text first so the month name is not
    # auto-converted into a date serial value.
    $ws.Range($col + "1").NumberFormat = "@"
    $ws.Range($col + "1").Value = $monthNames[$i]
}

# Numeric data rows 2-13; only Jan-Apr 2022 (PR:PU) have published values.
# May-Dec 2022 (PV:QC) are formatted but left blank (not yet available).
$rowValues = @{
    2  = @(278075,   212402,   293883,   347948.68)
    3  = @(255294,   197805,   244096,   306896)
    4  = @(195775,   136607,   170384,   239604)
    5  = @(58346,    59433,    70778,    65669)
    6  = @(1173,     1765,     2934,     1623)
    7  = @(22779,    14538,    49782,    40937.68)
    8  = @(2,        59,       5,        115)
    9  = @(301457,   318202,   481549,   343013)
    10 = @(79922,    93367,    94067,    85507)
    11 = @(65551,    28230,    55548,    37303)
    12 = @(149732,   188921,   318734,   207990)
    13 = @(-23382,   -105800,  -187666,  4935.679999999993)
}

$dataCols = @("PR","PS","PT","PU")
$blankCols = @("PV","PW","PX","PY","PZ","QA","QB","QC")

foreach ($row in 2..13) {
    $vals = $rowValues[$row]
    for ($i = 0; $i -lt $dataCols.Length; $i++) {
        $cell = $ws.Range($dataCols[$i] + $row)
        $cell.NumberFormat = "#,##0"
        $cell.Value = $vals[$i]
    }
    foreach ($col in $blankCols) {
        $ws.Range($col + $row).NumberFormat = "#,##0"
    }
}

# Move the frozen-pane scroll / active selection to reflect the newly
# entered data (mirrors the author ending up past the new columns).
$ws.Activate()
$ws.Range("PV14").Select()

# ---------------------------------------------------------------------
# Sheet "metadata": refresh "Latest Data" + "Availability" labels
# ---------------------------------------------------------------------
$meta = $wb.Worksheets.Item("metadata")

$meta.Range("B6").Value = "April 2022"

$availabilityRows = @(9, 10, 11, 14, 16, 17, 24, 25, 26, 27)
foreach ($r in $availabilityRows) {
    $meta.Range("B" + $r).Value = "1986-2022"
}

$meta.Range("B32").Value = "2000-2022 (Residual calculated for 1986-99)"

$meta.Activate()
$meta.Range("B33").Select()

# ---------------------------------------------------------------------
# Sheet "Annual": no data changes, just reflect end-of-edit selection.
# ---------------------------------------------------------------------
$annual = $wb.Worksheets.Item("Annual")
$annual.Activate()
$annual.Range("BF1").Select()

$ws.Activate()
